# Add new daily-tracker entries (rows 10-14) to the JAN-2021 sheet, matching
# rows 2-9's formatting (style indices + wrapped-text row height).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("JAN-2021")
$ws.Activate()

# Clone formatting (styles) from the last existing data row (row 9) down into
# the five new rows before writing any values into them.
$ws.Range("A9:G9").Copy()
$ws.Range("A10:G14").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Two-line wrapped task/comment text needs the taller row height (matches the
# existing wrapped rows); row 12 keeps the sheet's default row height.
$ws.Range("A10:G11").RowHeight = 28.8
$ws.Range("A13:G14").RowHeight = 28.8

# --- Row 10 -----------------------------------------------------------
$ws.Cells.Item(10, 1).Value = 9
$ws.Cells.Item(10, 2).Value = 44214
$ws.Cells.Item(10, 3).Value = "Selenium log Files(Sony, Samsung), Soniya"
$ws.Cells.Item(10, 4).Value = "Selenium log file Testing (QMVAR TO GSPN"
$ws.Cells.Item(10, 5).Value = 1
$ws.Cells.Item(10, 6).Value = "Completed"
$ws.Cells.Item(10, 7).Value = "Soniya setup Create(Modification)"

# --- Row 11 -----------------------------------------------------------
$ws.Cells.Item(11, 1).Value = 10
$ws.Cells.Item(11, 2).Value = 44215
$ws.Cells.Item(11, 3).Value = "Selenium log Files(Sony, Samsung), Git Hub"
$ws.Cells.Item(11, 4).Value = "Selenium log file Testing (QMVAR TO GSPN), Git Hub"
$ws.Cells.Item(11, 5).Value = 1
$ws.Cells.Item(11, 6).Value = "Completed"
$ws.Cells.Item(11, 7).Value = "Git Hub Upload Soniya Project Modules"

# --- Rows 12-14 ---------------------------------------------------------
# Task names (column C) for the QMVAR 2.0 work were entered together, then
# the shared "SETUP ADD USER" module note, then each day's comment.
$ws.Cells.Item(12, 3).Value = "QMVAR 2.0"
$ws.Cells.Item(13, 3).Value = "Selenium log Files(Sony, Samsung), QMVAR 2.0"
$ws.Cells.Item(12, 4).Value = "SETUP ADD USER"
$ws.Cells.Item(12, 7).Value = "Module, logic Create"
$ws.Cells.Item(13, 7).Value = "GridView Completed"
$ws.Cells.Item(14, 7).Value = "Create setup user completed"

$ws.Cells.Item(13, 4).Value = "SETUP ADD USER"
$ws.Cells.Item(14, 3).Value = "Selenium log Files(Sony, Samsung), QMVAR 2.0"
$ws.Cells.Item(14, 4).Value = "SETUP ADD USER"

$ws.Cells.Item(12, 1).Value = 11
$ws.Cells.Item(12, 2).Value = 44216
$ws.Cells.Item(12, 5).Value = 0.3
$ws.Cells.Item(12, 6).Value = "Completed"

$ws.Cells.Item(13, 1).Value = 12
$ws.Cells.Item(13, 2).Value = 44217
$ws.Cells.Item(13, 5).Value = 1
$ws.Cells.Item(13, 6).Value = "Completed"

$ws.Cells.Item(14, 1).Value = 13
$ws.Cells.Item(14, 2).Value = 44218
$ws.Cells.Item(14, 5).Value = 1
$ws.Cells.Item(14, 6).Value = "Completed"

# Update the visible selection / scroll position to match the edited sheet.
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C13").Select()

Write-Host "Added rows 10-14 to JAN-2021 sheet"
